# Update the "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" worksheets, which carry the
# same data table.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 69
    3  = 1052
    4  = 46
    6  = 2995
    8  = 2045
    10 = 105
    11 = 889
    13 = 31
    14 = 224
    16 = 93
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
